$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 55, shifting existing rows 55-104 down to 58-107
$ws.Range("A55:A57").EntireRow.Insert()

# Common constant values for the three new rows
$mercadoId = 1
$mercado = "Agrícola del Norte S.A. de Arica"
$region = "Arica y Parinacota"
$codreg = 15
$catId = 100112045
$categoria = "Zapallo"
$variedad = "Camote"
$unidad = "`$/kilo (volumen en unidades)"
$kgOUnidades = 1
$clasificacion = "Hortaliza"

# Row 55
$ws.Cells.Item(55,1).Value = $mercadoId
$ws.Cells.Item(55,2).Value = $mercado
$ws.Cells.Item(55,3).Value = $region
$ws.Cells.Item(55,4).Value = 45167
$ws.Cells.Item(55,5).Value = $codreg
$ws.Cells.Item(55,6).Value = $catId
$ws.Cells.Item(55,7).Value = $categoria
$ws.Cells.Item(55,8).Value = $variedad
$ws.Cells.Item(55,9).Value = "1a nueva(o)"
$ws.Cells.Item(55,10).Value = 850
$ws.Cells.Item(55,11).Value = 630
$ws.Cells.Item(55,12).Value = 650
$ws.Cells.Item(55,13).Value = 639
$ws.Cells.Item(55,14).Value = $unidad
$ws.Cells.Item(55,15).Value = "Perú"
$ws.Cells.Item(55,16).Value = 639
$ws.Cells.Item(55,17).Value = $kgOUnidades
$ws.Cells.Item(55,18).Value = $clasificacion

# Row 56
$ws.Cells.Item(56,1).Value = $mercadoId
$ws.Cells.Item(56,2).Value = $mercado
$ws.Cells.Item(56,3).Value = $region
$ws.Cells.Item(56,4).Value = 45167
$ws.Cells.Item(56,5).Value = $codreg
$ws.Cells.Item(56,6).Value = $catId
$ws.Cells.Item(56,7).Value = $categoria
$ws.Cells.Item(56,8).Value = $variedad
$ws.Cells.Item(56,9).Value = "2a nueva(o)"
$ws.Cells.Item(56,10).Value = 680
$ws.Cells.Item(56,11).Value = 600
$ws.Cells.Item(56,12).Value = 630
$ws.Cells.Item(56,13).Value = 617
$ws.Cells.Item(56,14).Value = $unidad
$ws.Cells.Item(56,15).Value = "Perú"
$ws.Cells.Item(56,16).Value = 617
$ws.Cells.Item(56,17).Value = $kgOUnidades
$ws.Cells.Item(56,18).Value = $clasificacion

# Row 57
$ws.Cells.Item(57,1).Value = $mercadoId
$ws.Cells.Item(57,2).Value = $mercado
$ws.Cells.Item(57,3).Value = $region
$ws.Cells.Item(57,4).Value = 45167
$ws.Cells.Item(57,5).Value = $codreg
$ws.Cells.Item(57,6).Value = $catId
$ws.Cells.Item(57,7).Value = $categoria
$ws.Cells.Item(57,8).Value = $variedad
$ws.Cells.Item(57,9).Value = "3a nueva (o)"
$ws.Cells.Item(57,10).Value = 300
$ws.Cells.Item(57,11).Value = 580
$ws.Cells.Item(57,12).Value = 600
$ws.Cells.Item(57,13).Value = 593
$ws.Cells.Item(57,14).Value = $unidad
$ws.Cells.Item(57,15).Value = "Perú"
$ws.Cells.Item(57,16).Value = 593
$ws.Cells.Item(57,17).Value = $kgOUnidades
$ws.Cells.Item(57,18).Value = $clasificacion
